# Update the "cryptos" worksheet with refreshed price / Volume(1h) figures,
# plus a few rows whose coin data was re-ordered (Stacks/FirstDigitalUSD/Hedera
# and ThetaToken/InjectiveProtocol), matching the GitHub Actions refresh commit.
#
# Note: Price values in column D are stored as plain text (e.g. "1.00", "6.80")
# so that trailing zeros / thousands-dot formatting are preserved exactly as in
# the source data. Assigning a numeric-looking string directly to .Value would
# make Excel auto-convert the cell to a real number (losing formatting), so a
# leading "'" is used to force text entry, and the style is then reset back to
# "Normal" to avoid leaving a visible Text number format on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.752.52"
$ws.Range("E2").Value = "  -3.84%  "
$ws.Range("D3").Value = "3.100.81"
$ws.Range("E3").Value = "  -4.76%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'608.95"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'145.07"
$ws.Range("E6").Value = "  -7.73%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.105.19"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = "  -4.36%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "  -7.53%  "
$ws.Range("D11").Value = "'5.25"
$ws.Range("E11").Value = "  -9.16%  "
$ws.Range("D12").Value = "'0.466"
$ws.Range("E12").Value = "  -5.97%  "
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  -8.14%  "
$ws.Range("D14").Value = "'35.03"
$ws.Range("E14").Value = "  -10.36%  "
$ws.Range("D15").Value = "3.645.01"
$ws.Range("E15").Value = "  -3.73%  "
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "63.781.92"
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("D18").Value = "3.120.65"
$ws.Range("E18").Value = "  -4.12%  "
$ws.Range("D19").Value = "'6.80"
$ws.Range("E19").Value = "  -8.54%  "
$ws.Range("D20").Value = "'473.35"
$ws.Range("E20").Value = "  -6.32%  "
$ws.Range("D21").Value = "'14.53"
$ws.Range("E21").Value = "  -5.99%  "
$ws.Range("D22").Value = "'0.695"
$ws.Range("E22").Value = "  -7.86%  "
$ws.Range("D23").Value = "'7.70"
$ws.Range("E23").Value = "  -4.78%  "
$ws.Range("D24").Value = "'13.51"
$ws.Range("E24").Value = "  -7.65%  "
$ws.Range("D25").Value = "'83.22"
$ws.Range("E25").Value = "  -4.95%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'2.78"
$ws.Range("E27").Value = "  -8.28%  "
$ws.Range("D28").Value = "'8.35"
$ws.Range("E28").Value = "  -9.08%  "
$ws.Range("D29").Value = "'2.14"
$ws.Range("E29").Value = "  -10.34%  "
$ws.Range("D30").Value = "'6.69"
$ws.Range("E30").Value = "  -4.91%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.73"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.110"
$ws.Range("E33").Value = "  -13.87%  "
$ws.Range("D34").Value = "'25.99"
$ws.Range("E34").Value = "  -7.11%  "
$ws.Range("D35").Value = "'1.10"
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").Value = "'5.89"
$ws.Range("E36").Value = "  -8.91%  "
$ws.Range("D37").Value = "'52.96"
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("D38").Value = "0.0₃0725"
$ws.Range("E38").Value = "  -7.23%  "
$ws.Range("D39").Value = "'457.34"
$ws.Range("E39").Value = "  -7.71%  "
$ws.Range("D40").Value = "'2.88"
$ws.Range("E40").Value = "  -14.60%  "
$ws.Range("D41").Value = "'0.0391"
$ws.Range("E41").Value = "  -7.29%  "
$ws.Range("D42").Value = "'0.119"
$ws.Range("E42").Value = "  -7.98%  "
$ws.Range("D43").Value = "'8.30"
$ws.Range("E43").Value = "  -6.17%  "
$ws.Range("D44").Value = "2.824.18"
$ws.Range("E44").Value = "  -5.83%  "
$ws.Range("D45").Value = "'0.263"
$ws.Range("E45").Value = "  -10.05%  "
$ws.Range("D46").Value = "'2.23"
$ws.Range("E46").Value = "  -11.19%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.36"
$ws.Range("E48").Value = "  -5.61%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'26.11"
$ws.Range("E49").Value = "  -9.93%  "
$ws.Range("E50").Value = "  -5.44%  "
$ws.Range("D51").Value = "'117.92"
$ws.Range("E51").Value = "  -2.15%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
